$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Cell="G2";Val=0.05317404096554723},
    @{Cell="H2";Val=25.77040556017616},
    @{Cell="I2";Val=37.55888262019705},
    @{Cell="G3";Val=0.06168516077963999},
    @{Cell="H3";Val=26.66266625458957},
    @{Cell="G4";Val=-0.001353261080915298},
    @{Cell="H4";Val=-171.0785356829792},
    @{Cell="G5";Val=0.01310433216054829},
    @{Cell="H5";Val=245.1823192069679},
    @{Cell="G6";Val=0.04081088027351821},
    @{Cell="H6";Val=17.71380245419932},
    @{Cell="G7";Val=0.07274923464881315},
    @{Cell="H7";Val=36.77509991538916},
    @{Cell="G8";Val=-0.01982648965361706},
    @{Cell="H8";Val=-5.343640884137939},
    @{Cell="G9";Val=0.007608050647025095},
    @{Cell="H9";Val=135.3404134866227},
    @{Cell="G10";Val=-0.0529067380102807},
    @{Cell="H10";Val=27.22730232242863},
    @{Cell="G11";Val=-0.09317683621115115},
    @{Cell="H11";Val=-1.267016068557512},
    @{Cell="G12";Val=-0.2159998724061861},
    @{Cell="H12";Val=11.64010277619457},
    @{Cell="G13";Val=-0.2812136524745084},
    @{Cell="H13";Val=-2.329011511091058},
    @{Cell="G14";Val=-0.04641075312079999},
    @{Cell="H14";Val=-25.10583310130388},
    @{Cell="G15";Val=0.01619364046126515},
    @{Cell="H15";Val=146.5717785687432},
    @{Cell="G16";Val=0.1203164502521125},
    @{Cell="H16";Val=-3.983568271934315},
    @{Cell="G17";Val=0.1531031276612776},
    @{Cell="H17";Val=9.164058087050252},
    @{Cell="G18";Val=0.1175712280463526},
    @{Cell="H18";Val=-5.741473231020283},
    @{Cell="G19";Val=0.11680660938743},
    @{Cell="H19";Val=-12.31646286613372},
    @{Cell="G20";Val=0.04714758823677372},
    @{Cell="H20";Val=37.31195673882298},
    @{Cell="G21";Val=0.08049770927079271},
    @{Cell="H21";Val=38.69383152486677},
    @{Cell="G22";Val=-0.08761473821967959},
    @{Cell="H22";Val=-9.726925723857612},
    @{Cell="G23";Val=-0.070245862208966},
    @{Cell="H23";Val=-12.30993879368396},
    @{Cell="G24";Val=0.1023225986924817},
    @{Cell="H24";Val=-13.3700011455833},
    @{Cell="G25";Val=0.1435464888595356},
    @{Cell="H25";Val=13.77192747281481},
    @{Cell="G26";Val=0.04904992289272603},
    @{Cell="H26";Val=-1.316159436647336},
    @{Cell="G27";Val=0.08881274147405324},
    @{Cell="H27";Val=2.464531260193028},
    @{Cell="G28";Val=-0.07739644553108113},
    @{Cell="H28";Val=-21.70393610355886},
    @{Cell="G29";Val=-0.09282872003794926},
    @{Cell="H29";Val=-30.42830274353686},
    @{Cell="G30";Val=0.08285637894780777},
    @{Cell="H30";Val=30.05621592922362},
    @{Cell="G31";Val=0.04796924637314676},
    @{Cell="H31";Val=-20.81729757123582},
    @{Cell="G32";Val=0.08341370026068774},
    @{Cell="H32";Val=-15.11264183818791},
    @{Cell="G33";Val=0.1144072422754222},
    @{Cell="H33";Val=39.03828431973768},
    @{Cell="G34";Val=-0.01232324823825591},
    @{Cell="H34";Val=-147.2966651079648},
    @{Cell="G35";Val=0.03902046196721032},
    @{Cell="H35";Val=448.0828095127285},
    @{Cell="G36";Val=0.004822201568103772},
    @{Cell="H36";Val=795.9751245909682},
    @{Cell="G37";Val=-0.001016843596958523},
    @{Cell="H37";Val=91.90041311289964},
    @{Cell="G38";Val=0.114065044069452},
    @{Cell="H38";Val=6.347147191508562},
    @{Cell="G39";Val=0.1160935791569639},
    @{Cell="H39";Val=35.52394224130772},
    @{Cell="G40";Val=0.02557158787815362},
    @{Cell="H40";Val=760.9216442120688},
    @{Cell="G41";Val=0.01874960816488109},
    @{Cell="H41";Val=25.02450614101122},
    @{Cell="G42";Val=0.09934409092929099},
    @{Cell="H42";Val=-1.57629426654976},
    @{Cell="G43";Val=0.1193266158267854},
    @{Cell="H43";Val=-0.6807530676840091},
    @{Cell="G44";Val=0.02222781845011079},
    @{Cell="H44";Val=-37.71554685500887},
    @{Cell="G45";Val=0.02122277206697491},
    @{Cell="H45";Val=29.64292166411527},
    @{Cell="G46";Val=0.05518558269945474},
    @{Cell="H46";Val=52.28432067353109},
    @{Cell="G47";Val=0.07455245384012942},
    @{Cell="H47";Val=47.80331073594636},
    @{Cell="G48";Val=0.07471951535936576},
    @{Cell="H48";Val=74.65460720818054},
    @{Cell="G49";Val=0.06468240872546509},
    @{Cell="H49";Val=-6.899771143180766},
    @{Cell="G50";Val=0.01239374283397024},
    @{Cell="H50";Val=-28.24666566963137},
    @{Cell="G51";Val=0.03325972195793012},
    @{Cell="H51";Val=70.82802039843183},
    @{Cell="G52";Val=-0.09270884550937053},
    @{Cell="H52";Val=10.44393579962433},
    @{Cell="G53";Val=-0.08855806027246961},
    @{Cell="H53";Val=4.111252924636609},
    @{Cell="G54";Val=0.08746364964312812},
    @{Cell="H54";Val=19.6122356621227},
    @{Cell="G55";Val=0.08513667762585654},
    @{Cell="H55";Val=37.42437456797359},
    @{Cell="G56";Val=0.02209789226496149},
    @{Cell="H56";Val=-36.84370367974066},
    @{Cell="G57";Val=0.006401134120666794},
    @{Cell="H57";Val=10.87060651407716},
    @{Cell="G58";Val=0.04511046429209091},
    @{Cell="H58";Val=80.36533212970697},
    @{Cell="G59";Val=0.02027707983176593},
    @{Cell="H59";Val=-14.36547027787497},
    @{Cell="G60";Val=0.02097743280312326},
    @{Cell="H60";Val=-35.3399687388994},
    @{Cell="G61";Val=0.01020469536638959},
    @{Cell="H61";Val=-19.38141197035454},
    @{Cell="G62";Val=0.07222091245121878},
    @{Cell="H62";Val=19.64509650421101},
    @{Cell="G63";Val=0.070079337841184},
    @{Cell="H63";Val=115.0365981449349},
    @{Cell="G64";Val=0.02288931602641685},
    @{Cell="H64";Val=-43.51968170717237},
    @{Cell="G65";Val=0.05719113234736888},
    @{Cell="H65";Val=2.014221374023387},
    @{Cell="G66";Val=0.09276411012212529},
    @{Cell="H66";Val=-0.8444252343586314},
    @{Cell="G67";Val=0.07621833771305943},
    @{Cell="H67";Val=-33.97960575773079},
    @{Cell="G68";Val=-0.02502512080103806},
    @{Cell="H68";Val=28.19290081275992},
    @{Cell="G69";Val=-0.01540015990267722},
    @{Cell="H69";Val=27.43253852616377},
    @{Cell="G70";Val=0.06577541254205625},
    @{Cell="H70";Val=-28.99603892419586},
    @{Cell="G71";Val=0.09587073694100312},
    @{Cell="H71";Val=5.111284034444527},
    @{Cell="G72";Val=-0.05946445453862801},
    @{Cell="H72";Val=-6.033107149116565},
    @{Cell="G73";Val=-0.07656199021610349},
    @{Cell="H73";Val=-3.795149205148341},
    @{Cell="G74";Val=0.1045170598540026},
    @{Cell="H74";Val=4.57268541683125},
    @{Cell="G75";Val=0.1475319747534058},
    @{Cell="H75";Val=51.46394438271133},
    @{Cell="G76";Val=0.02454732485074196},
    @{Cell="H76";Val=-4.001758182238276},
    @{Cell="G77";Val=0.002703945619253274},
    @{Cell="H77";Val=-80.83536802753596},
    @{Cell="G78";Val=0.09951696997094481},
    @{Cell="H78";Val=54.82523851422418},
    @{Cell="G79";Val=0.09360585370707446},
    @{Cell="H79";Val=22.01942320195528},
    @{Cell="G80";Val=-0.2118690307330025},
    @{Cell="H80";Val=-27.93374475528945},
    @{Cell="G81";Val=-0.1091125487676467},
    @{Cell="H81";Val=48.06439977228754},
    @{Cell="G82";Val=0.1327137998371498},
    @{Cell="H82";Val=15.70985682580607},
    @{Cell="G83";Val=0.1736071975087026},
    @{Cell="H83";Val=-2.458119233613961},
    @{Cell="G84";Val=0.07606511271635366},
    @{Cell="H84";Val=219.096537478981},
    @{Cell="G85";Val=0.06727430571410269},
    @{Cell="H85";Val=9.254207875097432}
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Val
}

Write-Host "Updated $($updates.Count) cells"
